$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")
$ws.Range("I7:M35").ClearContents()
$excel.CalculateFullRebuild()
